$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DM_Stat (column C) and P_Value (column D) values for rows 2-11
$ws.Range("C2").Value = 0.6786205614925603
$ws.Range("D2").Value = 0.5044549798162303

$ws.Range("C3").Value = 1.534416413199693
$ws.Range("D3").Value = 0.1391840534998763

$ws.Range("C4").Value = 0.98293925511605
$ws.Range("D4").Value = 0.3363244148837217

$ws.Range("C5").Value = 0.1926193927571476
$ws.Range("D5").Value = 0.8490240947374166

$ws.Range("C6").Value = 0.6391191143249563
$ws.Range("D6").Value = 0.5293443487944267

$ws.Range("C7").Value = 0.284448829224068
$ws.Range("D7").Value = 0.7787266743024626

$ws.Range("C8").Value = -0.316638919183646
$ws.Range("D8").Value = 0.7545025823282372

$ws.Range("C9").Value = -0.3984838508677284
$ws.Range("D9").Value = 0.6941146191634022

$ws.Range("C10").Value = -0.869263281372751
$ws.Range("D10").Value = 0.3940881583839038

$ws.Range("C11").Value = -0.5465133057309075
$ws.Range("D11").Value = 0.5902135200877194

$wb.Save()
